$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# New row 12: a "submit" item
# (order of writes matters for shared-string table ordering)
$ws.Range("C12").Value = "submit"
$ws.Range("C12").WrapText = $true

$ws.Range("D12").Value = "bla"
$ws.Range("D12").WrapText = $true

$ws.Range("F12").Value = "Go on!"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").WrapText = $true

# New row 13: another "note" item
$ws.Range("C13").Value = "note"
$ws.Range("C13").WrapText = $true

$ws.Range("F13").Value = "Good work, chap!"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").WrapText = $true

$ws.Range("D13").Value = "xx"
$ws.Range("D13").WrapText = $true

# The former "instruction" cell C2 becomes a "note" item instead.
$ws.Range("C2").Value = "note"

# Update selection from H16 to C3
$ws.Range("C3").Select()
